$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'298.91"
$ws.Range("E2").Value = "'2.94%"
$ws.Range("D3").Value = "'41.01"
$ws.Range("E3").Value = "'3.54%"
$ws.Range("D4").Value = "'5.047"
$ws.Range("E4").Value = "'0.37%"
$ws.Range("D5").Value = "'0.07475"
$ws.Range("E5").Value = "'1.82%"
$ws.Range("D6").Value = "'4.352"
$ws.Range("E6").Value = "'1.52%"
$ws.Range("D7").Value = "'1.588"
$ws.Range("E7").Value = "'2.30%"
$ws.Range("D8").Value = "'0.9391"
$ws.Range("E8").Value = "'3.04%"
$ws.Range("D10").Value = "'0.1209"
$ws.Range("E10").Value = "'1.46%"
$ws.Range("D11").Value = "'0.1817"
$ws.Range("E11").Value = "'4.30%"
$ws.Range("D12").Value = "'0.08778"
$ws.Range("E12").Value = "'0.96%"
$ws.Range("D13").Value = "'0.04230"
$ws.Range("E13").Value = "'1.84%"
$ws.Range("E14").Value = "'-0.31%"
$ws.Range("E15").Value = "'-1.23%"
$ws.Range("D16").Value = "'0.005968"
$ws.Range("E16").Value = "'1.44%"
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D17").Value = "'0.003875"
$ws.Range("E17").Value = "'-0.38%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.359"
$ws.Range("E18").Value = "'-1.01%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3307"
$ws.Range("E19").Value = "'0.76%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'7.905"
$ws.Range("E20").Value = "'4.44%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1372"
$ws.Range("E21").Value = "'1.47%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2951"
$ws.Range("E22").Value = "'2.34%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04005"
$ws.Range("E23").Value = "'4.23%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "'0.001261"
$ws.Range("E24").Value = "'-0.69%"
$ws.Range("D25").Value = "'0.0001223"
$ws.Range("E25").Value = "'-4.54%"
$ws.Range("D26").Value = "'0.0003703"
$ws.Range("E26").Value = "'-0.65%"
$ws.Range("D38").Value = "'0.02418"
$ws.Range("E38").Value = "'3.64%"
$ws.Range("D39").Value = "'0.05176"
$ws.Range("E39").Value = "'3.14%"
$ws.Range("D40").Value = "'0.006067"
$ws.Range("E40").Value = "'18.74%"
$ws.Range("D41").Value = "'0.007737"
$ws.Range("E41").Value = "'0.56%"
$ws.Range("D42").Value = "'0.1325"
$ws.Range("E42").Value = "'4.18%"
$ws.Range("D43").Value = "'0.007334"
$ws.Range("E43").Value = "'-0.52%"
$ws.Range("D44").Value = "'0.007160"
$ws.Range("E44").Value = "'2.86%"
$ws.Range("D45").Value = "'0.2975"
$ws.Range("E45").Value = "'-5.31%"
$ws.Range("D46").Value = "'0.00006245"
$ws.Range("E46").Value = "'-4.11%"
$ws.Range("D47").Value = "'0.00000000746"
$ws.Range("E47").Value = "'-0.65%"
$ws.Range("D48").Value = "'0.04682"
$ws.Range("E48").Value = "'-81.41%"
$ws.Range("D49").Value = "'0.004178"
$ws.Range("E49").Value = "'-0.64%"
$ws.Range("D50").Value = "'0.00002089"
$ws.Range("E50").Value = "'-0.65%"
$ws.Range("D51").Value = "'0.0001990"
$ws.Range("E51").Value = "'-0.65%"
